# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff", and the generated/handoff timestamps were refreshed.
# Touching the cells (rather than poking the shared-string table) updates
# every row that shows the status/date, exactly like the OOXML the CI job
# re-emits after running the report generator again.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language Status columns + the "Latest HO Xliff
#     Generate Date" column for the single tracked file -------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-25 22:57:11"

# --- zh-cn sheet: Status + Latest Handoff Datetime ------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-25 22:57:05"

# --- de-de sheet: Status + Latest Handoff Datetime ------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-25 22:57:11"

# --- Widen the now-longer "Status"/language columns so "Ready for handoff"
#     fits without truncating (mirrors the col width bump in the diff) ----
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
